# ROM Progress tracker update:
#  - Bugfix + basic block fetching + continue paper review
#  - Flip a batch of "Status" cells from their previous color to
#    Fonctionnel (green) / Problème (yellow), clear out several stale
#    "Notes" back to N/A, and log a new note about the Ninja demo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Colors used by the legend (F6:G9) on this sheet.
$green  = 5287936   # RGB(0,176,80)   -> "Fonctionnel"
$yellow = 65535      # RGB(255,255,0) -> "Problème"

# --- Status column (C): recolor rows that progressed / regressed ---
$ws.Range("C7").Interior.Color  = $green
$ws.Range("C10").Interior.Color = $green
$ws.Range("C12").Interior.Color = $green
$ws.Range("C13").Interior.Color = $green
$ws.Range("C15").Interior.Color = $green
$ws.Range("C23").Interior.Color = $green
$ws.Range("C24").Interior.Color = $green
$ws.Range("C27").Interior.Color = $green
$ws.Range("C28").Interior.Color = $green
$ws.Range("C29").Interior.Color = $green
$ws.Range("C30").Interior.Color = $green
$ws.Range("C35").Interior.Color = $green
$ws.Range("C26").Interior.Color = $yellow

# --- Notes column (D): clear resolved bugs back to N/A ... ---
$ws.Range("D7").Value  = "N/A"
$ws.Range("D10").Value = "N/A"
$ws.Range("D12").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("D15").Value = "N/A"
$ws.Range("D35").Value = "N/A"

# ... and log the newly found Ninja issue.
$ws.Range("D26").Value = "Kunais apparaîssent proches du joueur"

# --- Update the view / selection to where the review left off ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G22").Select()
